# Apply the edits described in the diff to the single worksheet workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row 1 (table 1 header) ---
$ws.Range("G1").Value = "table_header_position"

# --- Table 1 (rows 2-13): G-column descriptor text changes for every data row ---
$oldTable1Text = '"Duxford Range Part Number Description Dimensions Power Lumens Colour Temp. - Can be found on the center right position of the page"'
$newTable1Text = '"Part Number - Can be found on the top right position of the page"'

foreach ($r in 2..13) {
    $cell = $ws.Cells.Item($r, 7)  # column G
    if ($cell.Value() -eq $oldTable1Text) {
        $cell.Value = $newTable1Text
    }
}

# --- Table 1 Power/Lumens (D/E) cell value shuffle ---
# Row 3 gains 16W / 1600lm
$ws.Range("D3").Value = "16W"
$ws.Range("E3").Value = "1600lm"

# Row 11 gains 16W / 1600lm
$ws.Range("D11").Value = "16W"
$ws.Range("E11").Value = "1600lm"

# Row 12 loses its 16W / 1600lm values (becomes blank)
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""

# Row 13 gains 16W / 1600lm
$ws.Range("D13").Value = "16W"
$ws.Range("E13").Value = "1600lm"

# --- Header row 16 (table 2 header) ---
$ws.Range("G16").Value = "table_header_position"

# --- Table 2 (rows 17-28): G-column descriptor text changes for every data row ---
$oldTable2Text = '"Duxford Range - Can be found on the right side of the page"'
$newTable2Text = '"Part Number Description Dimensions Power Lumens Colour Temp. - Can be found on the right side of the page"'

foreach ($r in 17..28) {
    $cell = $ws.Cells.Item($r, 7)  # column G
    if ($cell.Value() -eq $oldTable2Text) {
        $cell.Value = $newTable2Text
    }
}
